$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their original text formatting
# (avoids Excel auto-converting numeric-looking strings into numbers)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "91.184.96"
$ws.Range("E2").Value = "  +4.26%  "
$ws.Range("D3").Value = "3.221.62"
$ws.Range("E3").Value = "  +2.01%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "218.87"
$ws.Range("E5").Value = "  +6.47%  "
$ws.Range("D6").Value = "650.76"
$ws.Range("E6").Value = "  +6.92%  "
$ws.Range("D7").Value = "0.400"
$ws.Range("E7").Value = "  +4.14%  "
$ws.Range("D8").Value = "0.702"
$ws.Range("E8").Value = "  +5.21%  "
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").Value = "3.217.46"
$ws.Range("E10").Value = "  +1.95%  "
$ws.Range("D11").Value = "0.581"
$ws.Range("E11").Value = "  +8.86%  "
$ws.Range("D12").Value = "0.181"
$ws.Range("E12").Value = "  +1.98%  "
$ws.Range("D13").Value = "0.0000259"
$ws.Range("E13").Value = "  +6.12%  "
$ws.Range("D14").Value = "5.45"
$ws.Range("E14").Value = "  +3.75%  "
$ws.Range("D15").Value = "33.65"
$ws.Range("E15").Value = "  +5.21%  "
$ws.Range("D16").Value = "90.722.65"
$ws.Range("E16").Value = "  +4.08%  "
$ws.Range("D17").Value = "3.816.30"
$ws.Range("E17").Value = "  +2.10%  "
$ws.Range("D18").Value = "3.220.32"
$ws.Range("E18").Value = "  +2.27%  "
$ws.Range("B19").Value = "PEPE"
$ws.Range("C19").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D19").Value = "0.0000229"
$ws.Range("E19").Value = "  +75.78%  "
$ws.Range("B20").Value = "SuiNetwork"
$ws.Range("C20").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D20").Value = "3.39"
$ws.Range("E20").Value = "  +13.14%  "
$ws.Range("D21").Value = "13.61"
$ws.Range("E21").Value = "  +1.71%  "
$ws.Range("D22").Value = "443.15"
$ws.Range("E22").Value = "  +6.93%  "
$ws.Range("D23").Value = "8.74"
$ws.Range("E23").Value = "  +3.07%  "
$ws.Range("D24").Value = "5.14"
$ws.Range("E24").Value = "  +0.99%  "
$ws.Range("D25").Value = "5.33"
$ws.Range("E25").Value = "  +3.09%  "
$ws.Range("D26").Value = "11.97"
$ws.Range("E26").Value = "  +1.26%  "
$ws.Range("D27").Value = "82.35"
$ws.Range("E27").Value = "  +12.16%  "
$ws.Range("D28").Value = "3.383.64"
$ws.Range("E28").Value = "  +1.79%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").Value = "0.162"
$ws.Range("E30").Value = "  +1.35%  "
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("B32").Value = "dogwifhat"
$ws.Range("C32").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D32").Value = "4.16"
$ws.Range("E32").Value = "  +39.03%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "550.20"
$ws.Range("E33").Value = "  +1.90%  "
$ws.Range("E34").Value = "  +3.72%  "
$ws.Range("D35").Value = "7.12"
$ws.Range("E35").Value = "  +5.80%  "
$ws.Range("E36").Value = "  +6.14%  "
$ws.Range("D37").Value = "1.31"
$ws.Range("E37").Value = "  +1.05%  "
$ws.Range("E38").Value = "  +3.55%  "
$ws.Range("D39").Value = "22.45"
$ws.Range("E39").Value = "  +3.01%  "
$ws.Range("D40").Value = "0.128"
$ws.Range("E40").Value = "  -2.92%  "
$ws.Range("E41").Value = "  +0.34%  "
$ws.Range("D42").Value = "1.96"
$ws.Range("E42").Value = "  +3.80%  "
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").Value = "0.378"
$ws.Range("E44").Value = "  +2.04%  "
$ws.Range("D45").Value = "45.30"
$ws.Range("E45").Value = "  +4.61%  "
$ws.Range("D46").Value = "146.72"
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("D47").Value = "174.95"
$ws.Range("E47").Value = "  +1.56%  "
$ws.Range("D48").Value = "0.769"
$ws.Range("E48").Value = "  +10.47%  "
$ws.Range("D49").Value = "0.126"
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("E50").Value = "  +1.45%  "
$ws.Range("D51").Value = "0.628"
$ws.Range("E51").Value = "  +7.26%  "
